$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$q1.Name = "2022-Q1"

# Header row (row 1) - matches the style used on the "2021-Q4" sheet
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("B1:H1").Style = $q1.Range("A1").Style

# Copy header style (bold/border/centered) from the "2021-Q4" sheet header
$hdrSrc = $afterSheet.Range("B1:D1")
$hdrSrc.Copy()
$q1.Range("B1:D1").PasteSpecial(-4122)
$afterSheet.Range("C1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Row 2 - 000593
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'000593"
$q1.Range("C2").Value = "易方达标普全球高端消费品指数增强(QDII)-美元现汇"
$q1.Range("D2").Value = "'1.93"
$q1.Range("E2").Value = "'92.46"
$q1.Range("F2").Value = "'5.19"
$q1.Range("G2").Value = "'0.1002"
$q1.Range("H2").Value = 9

# Row 3 - 005676
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'005676"
$q1.Range("C3").Value = "易方达标普全球高端消费品指数增强C(QDII) - 人民币"
$q1.Range("D3").Value = "'1.93"
$q1.Range("E3").Value = "'92.46"
$q1.Range("F3").Value = "'5.19"
$q1.Range("G3").Value = "'0.1002"
$q1.Range("H3").Value = 9

# Row 4 - 118002
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'118002"
$q1.Range("C4").Value = "易方达标普全球高端消费品指数增强A(QDII) - 人民币"
$q1.Range("D4").Value = "'1.93"
$q1.Range("E4").Value = "'92.46"
$q1.Range("F4").Value = "'5.19"
$q1.Range("G4").Value = "'0.1002"
$q1.Range("H4").Value = 9

# Match the "A" column style (bold/border/centered header style) used for
# the row-index column on the other sheets
$afterSheet.Range("A2").Copy()
$q1.Range("A2:A4").PasteSpecial(-4122)
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet - insert a new row for 2022-Q1
#    above the existing "2021-Q4" row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.3

# restore the row-index cell style to match the rest of column A
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0

# re-number the row-index column (0, 1, 2) after the insert
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# keep the originally-active sheet/tab selected, since adding a new
# worksheet would otherwise shift focus onto it
$wb.Worksheets.Item("2021-Q3").Activate()
